$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "70. Climbing Stairs"
$ws.Range("B2").Value = "recursive, DP. 2 ways to get to n step, take 2 step at n - 2, take 1 step at n - 1."

$ws.Range("B2").Select()

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
